$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.958.56'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '3.191.65'
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.32'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.74'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.54%  '

$ws.Range("D8").Value = '3.146.44'
$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.82'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.04'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.98%  '

$ws.Range("D15").Value = '3.716.02'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("E16").Value = '  -1.59%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '63.918.42'
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.157.28'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.08'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.74'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.33'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.734'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.44'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.94'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.01'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.23'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.23%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.35'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +8.21%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.68'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.03%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.26%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.12'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.99'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.51%  '

$ws.Range("E34").Value = '  +0.49%  '

$ws.Range("D35").Value = '0.0₃0860'
$ws.Range("E35").Value = '  -1.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.32'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.35'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.02'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '51.35'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '437.33'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.89'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.02%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0372'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.29%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.918.97'
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.283'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.108'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '37.84'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +11.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.80'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.84%  '

$ws.Range("E50").Value = '  -1.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.17'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.92%  '
